# Updated cryptos list price/volume figures (GitHub Actions scrape refresh).
# Column D = Price (text, may look numeric -> force text via leading apostrophe + restore
# default style so Excel does not coerce it to a Number cell or leave a quote-prefix style).
# Column E = Volume(1h) (already non-numeric text like "  -3.71%  ", safe to set directly).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '24.942.51'
$ws.Range("E2").Value = '  -3.71%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.636.87'
# Row 4: TetherUSD
$ws.Range("D4").Value = '''0.9971'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.36%  '
# Row 5: BNB
$ws.Range("D5").Value = '''236.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.49%  '
# Row 6: USDC
$ws.Range("E6").Value = '  -0.01%  '
# Row 7: XRP
$ws.Range("D7").Value = '''0.4700'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.12%  '
# Row 8: Cardano
$ws.Range("D8").Value = '''0.2560'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.05%  '
# Row 9: Dogecoin
$ws.Range("D9").Value = '''0.06010'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.70%  '
# Row 10: TRON
$ws.Range("D10").Value = '''0.07094'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.06%  '
# Row 11: WrappedEther
$ws.Range("D11").Value = '1.636.07'
$ws.Range("E11").Value = '  -6.14%  '
# Row 12: Solana
$ws.Range("D12").Value = '''14.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.58%  '
# Row 13: Polygon
$ws.Range("E13").Value = '  -5.36%  '
# Row 14: Polkadot
$ws.Range("D14").Value = '''4.405'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.64%  '
# Row 15: Litecoin
$ws.Range("D15").Value = '''72.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.89%  '
# Row 16: Dai
$ws.Range("E16").Value = '  -0.01%  '
# Row 17: BinanceUSD
$ws.Range("D17").Value = '''0.9974'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.29%  '
# Row 18: WrappedBTC
$ws.Range("D18").Value = '24.934.43'
$ws.Range("E18").Value = '  -3.81%  '
# Row 19: ShibaInu
$ws.Range("D19").Value = '''0.000006582'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.11%  '
# Row 20: Avalanche
$ws.Range("D20").Value = '''11.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.24%  '
# Row 21: Uniswap
$ws.Range("D21").Value = '''4.404'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.56%  '
# Row 22: WrappedliquidstakedEther2.0
$ws.Range("D22").Value = '1.843.93'
$ws.Range("E22").Value = '  -6.53%  '
# Row 23: Cosmos
$ws.Range("D23").Value = '''8.598'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.37%  '
# Row 24: Chainlink
$ws.Range("D24").Value = '''5.262'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.05%  '
# Row 25: Monero
$ws.Range("D25").Value = '''132.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.46%  '
# Row 26: EthereumClassic
$ws.Range("D26").Value = '''14.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.58%  '
# Row 27: Toncoin
$ws.Range("D27").Value = '''1.365'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.69%  '
# Row 28: BitcoinCash
$ws.Range("D28").Value = '''102.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.89%  '
# Row 29: LidoDAOToken
$ws.Range("D29").Value = '''1.655'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.01%  '
# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = '''3.734'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.31%  '
# Row 31: Stellar
$ws.Range("D31").Value = '''0.07722'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.13%  '
# Row 32: Filecoin
$ws.Range("D32").Value = '''3.551'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.97%  '
# Row 33: Hedera
$ws.Range("D33").Value = '''0.04343'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.99%  '
# Row 34: Frax
$ws.Range("D34").Value = '''0.9989'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.06%  '
# Row 35: HuobiToken
$ws.Range("D35").Value = '''2.600'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.05%  '
# Row 36: ARBITRUM
$ws.Range("D36").Value = '''0.9213'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.11%  '
# Row 37: ImmutableX
$ws.Range("D37").Value = '''0.5802'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.23%  '
# Row 38: MXToken
$ws.Range("D38").Value = '''2.540'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.86%  '
# Row 39: VeChain
$ws.Range("E39").Value = '  -3.09%  '
# Row 40: PaxDollar
$ws.Range("D40").Value = '''0.9976'
$ws.Range("D40").Style = "Normal"
# Row 41: TrustWalletToken
$ws.Range("D41").Value = '''0.8257'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.54%  '
# Row 42: RenderToken
$ws.Range("D42").Value = '''1.793'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.05%  '
# Row 43: Quant
$ws.Range("D43").Value = '''97.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.49%  '
# Row 44: TheSandbox
$ws.Range("D44").Value = '''0.3719'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.55%  '
# Row 45: FraxShare
$ws.Range("D45").Value = '''4.727'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.99%  '
# Row 46: Algorand
$ws.Range("D46").Value = '''0.1114'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.27%  '
# Row 47: Cronos
$ws.Range("D47").Value = '''0.05226'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
# Row 48: Aptos
$ws.Range("D48").Value = '''6.082'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.89%  '
# Row 49: Elrond
$ws.Range("D49").Value = '''29.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.10%  '
# Row 50: TrueUSD
$ws.Range("D50").Value = '''0.9983'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
# Row 51: USDD
$ws.Range("D51").Value = '''0.9992'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.53%  '
